$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Single-occurrence line-number bumps
Replace-Text "M2DocEvaluator.caseQuery(M2DocEvaluator.java:587)" "M2DocEvaluator.caseQuery(M2DocEvaluator.java:591)"
Replace-Text "M2DocEvaluator.caseBlock(M2DocEvaluator.java:1467)" "M2DocEvaluator.caseBlock(M2DocEvaluator.java:1556)"
Replace-Text "M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:297)" "M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:301)"
Replace-Text "M2DocEvaluator.generate(M2DocEvaluator.java:282)" "M2DocEvaluator.generate(M2DocEvaluator.java:286)"
Replace-Text "M2DocUtils.generate(M2DocUtils.java:845)" "M2DocUtils.generate(M2DocUtils.java:853)"
Replace-Text "AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:514)" "AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:518)"
Replace-Text "AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:421)" "AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:414)"

# Three occurrences of the same old text -> same new text
Replace-Text "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1242)" "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1331)"

# Replace the tail of the JUnit/Eclipse launcher frames with the Maven/Tycho frames
$oldTail = "`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)`n" + `
    "`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)`n" + `
    "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)`n" + `
    "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)`n" + `
    "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)`n" + `
    "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)"

$newTail = "`tat org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:365)`n" + `
    "`tat org.apache.maven.surefire.junit4.JUnit4Provider.executeWithRerun(JUnit4Provider.java:273)`n" + `
    "`tat org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:238)`n" + `
    "`tat org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:159)`n" + `
    "`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n" + `
    "`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n" + `
    "`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n" + `
    "`tat java.lang.reflect.Method.invoke(Method.java:498)`n" + `
    "`tat org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:206)`n" + `
    "`tat org.apache.maven.surefire.booter.ProviderFactory`$ProviderProxy.invoke(ProviderFactory.java:161)`n" + `
    "`tat org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:84)`n" + `
    "`tat org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:113)`n" + `
    "`tat org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)`n" + `
    "`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n" + `
    "`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n" + `
    "`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n" + `
    "`tat java.lang.reflect.Method.invoke(Method.java:498)`n" + `
    "`tat org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:593)`n" + `
    "`tat org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:205)`n" + `
    "`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:137)`n" + `
    "`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:107)`n" + `
    "`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:401)`n" + `
    "`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:255)`n" + `
    "`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n" + `
    "`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n" + `
    "`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n" + `
    "`tat java.lang.reflect.Method.invoke(Method.java:498)`n" + `
    "`tat org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:657)`n" + `
    "`tat org.eclipse.equinox.launcher.Main.basicRun(Main.java:594)`n" + `
    "`tat org.eclipse.equinox.launcher.Main.run(Main.java:1447)`n" + `
    "`tat org.eclipse.equinox.launcher.Main.main(Main.java:1420)"

Replace-Text $oldTail $newTail
